$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 683, shifting existing rows 683..751 down to 684..752.
$ws.Rows(683).Insert()

# Populate the newly inserted row 683 with its data.
$row = 683
$ws.Cells.Item($row, 1).Value = 9
$ws.Cells.Item($row, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item($row, 3).Value = "Metropolitana"
$ws.Cells.Item($row, 4).Value = 45194
$ws.Cells.Item($row, 5).Value = 13
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108002
$ws.Cells.Item($row, 10).Value = "Mango"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 550
$ws.Cells.Item($row, 14).Value = 10500
$ws.Cells.Item($row, 15).Value = 11500
$ws.Cells.Item($row, 16).Value = 11045
$ws.Cells.Item($row, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item($row, 18).Value = "Brasil"
$ws.Cells.Item($row, 19).Value = 2761
$ws.Cells.Item($row, 20).Value = 4
